$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 currently holds the "Americana (o)" / 2021-12-29 record. The edit
# duplicates that record down into a new row 22 (unchanged), and then turns
# row 21 into a new weekly record dated 2023-08-16 with variety
# "Sin especificar".

# 1) Copy row 21's values (and D's date format) cell by cell into the new
#    row 22, so row 22 starts out identical to the current row 21.
for ($col = 1; $col -le 18; $col++) {
    $src = $ws.Cells.Item(21, $col)
    $dst = $ws.Cells.Item(22, $col)
    $dst.Value = $src.Value2
}
$ws.Cells.Item(22, 4).NumberFormat = $ws.Cells.Item(21, 4).NumberFormat

# 2) Update row 21 in place with the new weekly values.
$ws.Range("D21").Value = 45154
$ws.Range("H21").Value = "Sin especificar"
